# The extra R-SQL ETL code is removed from HTN pdf and word tsd
#
# Renames a handful of code_set identifiers on Sheet1:
#   - Schizophrenia/Bipolar row: drop the "_c_icd" infix
#       schizophrenia_c_icd_01_base -> schizophrenia_01_base
#       schizophrenia_c_icd_01      -> schizophrenia_01
#       bipolar_c_icd_01_base       -> bipolar_01_base
#       bipolar_c_icd_01            -> bipolar_01
#   - Covid Vaccination Pediatrics tsd file name: "vaccination" -> "vaccinations"
#       covid_vaccination_ped_covid_02_tsd_01 -> covid_vaccinations_ped_covid_02_tsd_01

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Keep shared-string append order matching the authored edit: the tsd file
# name first, then the schizophrenia/bipolar base & mapped code_sets.
$ws.Range("D17").Value = "covid_vaccinations_ped_covid_02_tsd_01"

$ws.Range("E6").Value = "schizophrenia_01_base"
$ws.Range("F6").Value = "schizophrenia_01"
$ws.Range("E7").Value = "bipolar_01_base"
$ws.Range("F7").Value = "bipolar_01"

# Match the author's final selection/view state in the saved file.
$ws.Range("F7").Select()
